$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.079.27'
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").Value = '1.787.63'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''221.99'
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '''32.25'
$ws.Range("E8").Value = '  -0.35%  '
$ws.Range("D9").Value = '''0.285'
$ws.Range("E9").Value = '  -0.44%  '
$ws.Range("D10").Value = '''0.0714'
$ws.Range("E10").Value = '  -1.22%  '
$ws.Range("D11").Value = '''0.0930'
$ws.Range("E11").Value = '  -0.22%  '
$ws.Range("D12").Value = '2.045.91'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '1.798.06'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("D14").Value = '''10.84'
$ws.Range("E14").Value = '  -1.65%  '
$ws.Range("D15").Value = '''0.627'
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("D16").Value = '34.054.26'
$ws.Range("E16").Value = '  -1.97%  '
$ws.Range("D17").Value = '''4.17'
$ws.Range("E17").Value = '  -3.26%  '
$ws.Range("D18").Value = '''68.06'
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D19").Value = '''243.90'
$ws.Range("E19").Value = '  -4.10%  '
$ws.Range("D20").Value = '0.0₃0784'
$ws.Range("E20").Value = '  -3.62%  '
$ws.Range("D21").Value = '''0.999'
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").Value = '''4.08'
$ws.Range("E23").Value = '  -3.54%  '
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("D25").Value = '''159.04'
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("D26").Value = '''16.34'
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("D27").Value = '''7.05'
$ws.Range("E27").Value = '  -0.94%  '
$ws.Range("D28").Value = '''0.112'
$ws.Range("E28").Value = '  -1.73%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  -2.28%  '
$ws.Range("E31").Value = '  +0.39%  '
$ws.Range("E32").Value = '  -3.64%  '
$ws.Range("E33").Value = '  -3.36%  '
$ws.Range("E34").Value = '  -3.48%  '
$ws.Range("D35").Value = '1.396.11'
$ws.Range("E35").Value = '  -3.01%  '
$ws.Range("D36").Value = '''0.645'
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("E38").Value = '  -3.37%  '
$ws.Range("D39").Value = '''79.56'
$ws.Range("E39").Value = '  -6.22%  '
$ws.Range("B40").Value = 'HuobiToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D40").Value = '''2.35'
$ws.Range("E40").Value = '  +1.05%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '''0.920'
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("D42").Value = '''2.71'
$ws.Range("E42").Value = '  -3.17%  '
$ws.Range("D43").Value = '''2.17'
$ws.Range("E43").Value = '  +2.21%  '
$ws.Range("D44").Value = '''0.0495'
$ws.Range("E44").Value = '  +0.23%  '
$ws.Range("E45").Value = '  -0.95%  '
$ws.Range("D46").Value = '''107.39'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("D48").Value = '1.944.27'
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("D50").Value = '''1.00'
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("E51").Value = '  +1.20%  '
